$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New local extremums: Anapa (rows 36-37) and Gelendzhik (rows 38-39)
$data = @(
    @(3703000, "город-курорт Анапа", "female", 2021, 0.06076, 0.06207, 0.05267, 0.04654, 0.03842, 0.0566, 0.0803, 0.08844, 0.0789, 0.0705, 0.1244, 0.07605, 0.0876, 0.0768),
    @(3703000, "город-курорт Анапа", "male", 2021, 0.07184, 0.0743, 0.0647, 0.0526, 0.04132, 0.0693, 0.08954, 0.0841, 0.07654, 0.0698, 0.11053, 0.0651, 0.0704, 0.05997),
    @(3708000, "город-курорт Геленджик", "female", 2021, 0.05774, 0.0602, 0.05148, 0.04257, 0.03543, 0.04437, 0.08014, 0.0925, 0.0836, 0.0747, 0.131, 0.0783, 0.0889, 0.07904),
    @(3708000, "город-курорт Геленджик", "male", 2021, 0.068, 0.07275, 0.06183, 0.05066, 0.0404, 0.05322, 0.08075, 0.09174, 0.08563, 0.07135, 0.1189, 0.0689, 0.073, 0.0628)
)

$startRow = 36
$endRow = 39

# Match the style of the existing data rows (s="1", center-aligned) for the
# newly appended rows before writing values into them.
$fmtSrc = $ws.Range("A2:R2")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A" + $startRow + ":R" + $endRow)
$fmtDst.PasteSpecial(-4122)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("T31").Select() | Out-Null
